$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table 1: rows 2-15, columns A-E and G (total change)
$ws.Range("A2").Value = 1.0
$ws.Range("B2").Value = "Lundo’s Legends"
$ws.Range("C2").Value = 133.5
$ws.Range("D2").Value = 60.0
$ws.Range("E2").Value = 73.5
$ws.Range("G2").Value = 0.10999999940395355

$ws.Range("A3").Value = 2.0
$ws.Range("B3").Value = "EL Onće"
$ws.Range("C3").Value = 132.0
$ws.Range("D3").Value = 65.5
$ws.Range("E3").Value = 66.5
$ws.Range("G3").Value = 0.10999999940395355

$ws.Range("A4").Value = 3.0
$ws.Range("B4").Value = "Samsquanches"
$ws.Range("C4").Value = 105.5
$ws.Range("D4").Value = 59.0
$ws.Range("E4").Value = 46.5
$ws.Range("G4").Value = 0.10999999940395355

$ws.Range("A5").Value = 4.0
$ws.Range("B5").Value = "GOD WILLS IT"
$ws.Range("C5").Value = 100.5
$ws.Range("D5").Value = 49.5
$ws.Range("E5").Value = 51.0
$ws.Range("G5").Value = 0.10999999940395355

$ws.Range("A6").Value = 5.0
$ws.Range("B6").Value = "Swampnuts"
$ws.Range("C6").Value = 95.0
$ws.Range("D6").Value = 39.0
$ws.Range("E6").Value = 56.0
$ws.Range("G6").Value = 0.10999999940395355

$ws.Range("A7").Value = 6.0
$ws.Range("B7").Value = "Splitfinger Skadoosh"
$ws.Range("C7").Value = 94.5
$ws.Range("D7").Value = 39.0
$ws.Range("E7").Value = 55.5
$ws.Range("G7").Value = 0.10999999940395355

$ws.Range("A8").Value = 7.0
$ws.Range("B8").Value = "Epic7"
$ws.Range("C8").Value = 89.0
$ws.Range("D8").Value = 40.0
$ws.Range("E8").Value = 49.0
$ws.Range("G8").Value = 0.10999999940395355

$ws.Range("A9").Value = 8.0
$ws.Range("B9").Value = "rainmaker"
$ws.Range("C9").Value = 81.5
$ws.Range("D9").Value = 23.0
$ws.Range("E9").Value = 58.5
$ws.Range("G9").Value = 0.10999999940395355

$ws.Range("A10").Value = 9.0
$ws.Range("B10").Value = "confusion"
$ws.Range("C10").Value = 78.0
$ws.Range("D10").Value = 52.0
$ws.Range("E10").Value = 26.0
$ws.Range("G10").Value = 0.10999999940395355

$ws.Range("A11").Value = 10.0
$ws.Range("B11").Value = "Mac"
$ws.Range("C11").Value = 76.5
$ws.Range("D11").Value = 44.0
$ws.Range("E11").Value = 32.5
$ws.Range("G11").Value = 0.10999999940395355

$ws.Range("A12").Value = 11.0
$ws.Range("B12").Value = "MillerTime"
$ws.Range("C12").Value = 72.5
$ws.Range("D12").Value = 38.0
$ws.Range("E12").Value = 34.5
$ws.Range("G12").Value = 0.10999999940395355

$ws.Range("A13").Value = 12.0
$ws.Range("B13").Value = "SmokeWalkers"
$ws.Range("C13").Value = 71.5
$ws.Range("D13").Value = 36.0
$ws.Range("E13").Value = 35.5
$ws.Range("G13").Value = 0.10999999940395355

$ws.Range("A14").Value = 13.0
$ws.Range("B14").Value = "Corbin Copy"
$ws.Range("C14").Value = 67.0
$ws.Range("D14").Value = 44.5
$ws.Range("E14").Value = 22.5
$ws.Range("G14").Value = 0.10999999940395355

$ws.Range("A15").Value = 14.0
$ws.Range("B15").Value = "DJ's Quality Team"
$ws.Range("C15").Value = 63.0
$ws.Range("D15").Value = 40.5
$ws.Range("E15").Value = 22.5
$ws.Range("G15").Value = 0.10999999940395355

# Table 2: rows 18-31, columns A-C and E-G
$ws.Range("A18").Value = 1.0
$ws.Range("B18").Value = "EL Onće"
$ws.Range("C18").Value = 65.5
$ws.Range("E18").Value = 1.0
$ws.Range("F18").Value = "Lundo’s Legends"
$ws.Range("G18").Value = 73.5

$ws.Range("A19").Value = 2.0
$ws.Range("B19").Value = "Lundo’s Legends"
$ws.Range("C19").Value = 60.0
$ws.Range("E19").Value = 2.0
$ws.Range("F19").Value = "EL Onće"
$ws.Range("G19").Value = 66.5

$ws.Range("A20").Value = 3.0
$ws.Range("B20").Value = "Samsquanches"
$ws.Range("C20").Value = 59.0
$ws.Range("E20").Value = 3.0
$ws.Range("F20").Value = "rainmaker"
$ws.Range("G20").Value = 58.5

$ws.Range("A21").Value = 4.0
$ws.Range("B21").Value = "confusion"
$ws.Range("C21").Value = 52.0
$ws.Range("E21").Value = 4.0
$ws.Range("F21").Value = "Swampnuts"
$ws.Range("G21").Value = 56.0

$ws.Range("A22").Value = 5.0
$ws.Range("B22").Value = "GOD WILLS IT"
$ws.Range("C22").Value = 49.5
$ws.Range("E22").Value = 5.0
$ws.Range("F22").Value = "Splitfinger Skadoosh"
$ws.Range("G22").Value = 55.5

$ws.Range("A23").Value = 6.0
$ws.Range("B23").Value = "Corbin Copy"
$ws.Range("C23").Value = 44.5
$ws.Range("E23").Value = 6.0
$ws.Range("F23").Value = "GOD WILLS IT"
$ws.Range("G23").Value = 51.0

$ws.Range("A24").Value = 7.0
$ws.Range("B24").Value = "Mac"
$ws.Range("C24").Value = 44.0
$ws.Range("E24").Value = 7.0
$ws.Range("F24").Value = "Epic7"
$ws.Range("G24").Value = 49.0

$ws.Range("A25").Value = 8.0
$ws.Range("B25").Value = "DJ's Quality Team"
$ws.Range("C25").Value = 40.5
$ws.Range("E25").Value = 8.0
$ws.Range("F25").Value = "Samsquanches"
$ws.Range("G25").Value = 46.5

$ws.Range("A26").Value = 9.0
$ws.Range("B26").Value = "Epic7"
$ws.Range("C26").Value = 40.0
$ws.Range("E26").Value = 9.0
$ws.Range("F26").Value = "SmokeWalkers"
$ws.Range("G26").Value = 35.5

$ws.Range("A27").Value = 10.5
$ws.Range("B27").Value = "Swampnuts"
$ws.Range("C27").Value = 39.0
$ws.Range("E27").Value = 10.0
$ws.Range("F27").Value = "MillerTime"
$ws.Range("G27").Value = 34.5

$ws.Range("A28").Value = 10.5
$ws.Range("B28").Value = "Splitfinger Skadoosh"
$ws.Range("C28").Value = 39.0
$ws.Range("E28").Value = 11.0
$ws.Range("F28").Value = "Mac"
$ws.Range("G28").Value = 32.5

$ws.Range("A29").Value = 12.0
$ws.Range("B29").Value = "MillerTime"
$ws.Range("C29").Value = 38.0
$ws.Range("E29").Value = 12.0
$ws.Range("F29").Value = "confusion"
$ws.Range("G29").Value = 26.0

$ws.Range("A30").Value = 13.0
$ws.Range("B30").Value = "SmokeWalkers"
$ws.Range("C30").Value = 36.0
$ws.Range("E30").Value = 13.5
$ws.Range("F30").Value = "Corbin Copy"
$ws.Range("G30").Value = 22.5

$ws.Range("A31").Value = 14.0
$ws.Range("B31").Value = "rainmaker"
$ws.Range("C31").Value = 23.0
$ws.Range("E31").Value = 13.5
$ws.Range("F31").Value = "DJ's Quality Team"
$ws.Range("G31").Value = 22.5
